$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.155.67'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.667.77'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5222'
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2605'
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06341'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.07'
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07545'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '1.671.66'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.425'
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5422'
$ws.Range("E14").Value = '  -5.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000007996'
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.46'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '26.185.27'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.729'
$ws.Range("E19").Value = '  -3.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.33'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.25'
$ws.Range("E21").Value = '  -3.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.226'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.58'
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1236'
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.425'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.76'
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06282'
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.362'
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.493'
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.410'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.639'
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9985'
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.393'
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.757'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5952'
$ws.Range("E37").Value = '  -2.59%  '
$ws.Range("D38").Value = '1.109.36'
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.060'
$ws.Range("E39").Value = '  -2.04%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01607'
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8562'
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.66'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").Value = '1.815.79'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.38'
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.055'
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05240'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4234'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.905'
$ws.Range("E51").Value = '  -1.48%  '
